$wb = $excel.ActiveWorkbook

# Rename the existing sheet from "Sheet1" to "Estimates"
$wsEstimates = $wb.Worksheets.Item("Sheet1")
$wsEstimates.Name = "Estimates"

# Add a new worksheet named "Effort" after the Estimates sheet
$wsEffort = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsEstimates)
$wsEffort.Name = "Effort"

# Populate the Effort sheet with header + data rows
$wsEffort.Range("A1").Value = "year"
$wsEffort.Range("B1").Value = "n_data_tot"
$wsEffort.Range("A2").Value = 1997
$wsEffort.Range("B2").Value = 3092
$wsEffort.Range("A3").Value = 1998
$wsEffort.Range("B3").Value = 2818

# Restore the selection on the Effort sheet to match the target state
$wsEffort.Range("B2").Select()

# Restore the selection on the Estimates sheet to match the target state
$wsEstimates.Activate()
$wsEstimates.Range("D35").Select()
